# Update staff data on Sheet1 (library, IQAC, Staff excel update)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add "Mr." prefix to the three staff names in column A (rows 8-10)
$ws.Range("A8").Value = "Mr. MANOHARAN P A"
$ws.Range("A9").Value = "Mr. VENKATESAN A"
$ws.Range("A10").Value = "Mr.GOKUL RAJ J"

# Move the active selection to A11 (clears the old D12 selection / B6 scroll position)
$ws.Range("A11").Select()
